$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.195798516273499
$ws.Range("B1").Value = 2.479881048202515
$ws.Range("C1").Value = 4.013492107391357
$ws.Range("D1").Value = 2.11854100227356
$ws.Range("E1").Value = 1.185153484344482
